$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("M2").Value = 2.565830333333333
$ws.Range("N2").Value = 7.697490999999999
$ws.Range("O2").Value = 0.0934185609347503
$ws.Range("P2").Value = 0.0934185609347503
$ws.Range("Q2").Value = 0.8035641779629998
$ws.Range("R2").Value = 7.232077601666999
$ws.Range("S2").Value = 0.002529274463938871
$ws.Range("T2").Value = 0.002529274463938871
$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("O3").Value = 0.3847798091300315
$ws.Range("P3").Value = 0.3847798091300315
$ws.Range("Q3").Value = 3.309784136327
$ws.Range("R3").Value = 29.788057226943
$ws.Range("S3").Value = 0.01041777710696719
$ws.Range("T3").Value = 0.01041777710696719
$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 13.68376133333333
$ws.Range("N4").Value = 41.051284
$ws.Range("O4").Value = 0.4982080363333638
$ws.Range("P4").Value = 0.4982080363333638
$ws.Range("Q4").Value = 4.285466690612
$ws.Range("R4").Value = 38.569200215508
$ws.Range("S4").Value = 0.01348880620102087
$ws.Range("T4").Value = 0.01348880620102087
$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6480206666666667
$ws.Range("N5").Value = 1.944062
$ws.Range("O5").Value = 0.0235935936018545
$ws.Range("P5").Value = 0.0235935936018545
$ws.Range("Q5").Value = 0.202946464366
$ws.Range("R5").Value = 1.826518179294
$ws.Range("S5").Value = 0.0006387881938301623
$ws.Range("T5").Value = 0.0006387881938301623
$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("M6").Value = 2.565830333333333
$ws.Range("N6").Value = 7.697490999999999
$ws.Range("O6").Value = 0.0934185609347503
$ws.Range("P6").Value = 0.0934185609347503
$ws.Range("Q6").Value = 20.72467687290044
$ws.Range("R6").Value = 186.522091856104
$ws.Range("S6").Value = 0.06523237026430162
$ws.Range("T6").Value = 0.06523237026430163
$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("O7").Value = 0.3847798091300315
$ws.Range("P7").Value = 0.3847798091300315
$ws.Range("S7").Value = 0.2686842820981699
$ws.Range("T7").Value = 0.2686842820981699
$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 13.68376133333333
$ws.Range("N8").Value = 41.051284
$ws.Range("O8").Value = 0.4982080363333638
$ws.Range("P8").Value = 0.4982080363333638
$ws.Range("Q8").Value = 110.5262216113885
$ws.Range("R8").Value = 994.7359945024961
$ws.Range("S8").Value = 0.3478890144480846
$ws.Range("T8").Value = 0.3478890144480847
$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6480206666666667
$ws.Range("N9").Value = 1.944062
$ws.Range("O9").Value = 0.0235935936018545
$ws.Range("P9").Value = 0.0235935936018545
$ws.Range("Q9").Value = 5.234180432414223
$ws.Range("R9").Value = 47.107623891728
$ws.Range("S9").Value = 0.01647494907116602
$ws.Range("T9").Value = 0.01647494907116602
$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("M10").Value = 2.565830333333333
$ws.Range("N10").Value = 7.697490999999999
$ws.Range("O10").Value = 0.0934185609347503
$ws.Range("P10").Value = 0.0934185609347503
$ws.Range("Q10").Value = 7.435123729818555
$ws.Range("R10").Value = 66.916113568367
$ws.Range("S10").Value = 0.02340257206801707
$ws.Range("T10").Value = 0.02340257206801707
$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("O11").Value = 0.3847798091300315
$ws.Range("P11").Value = 0.3847798091300315
$ws.Range("Q11").Value = 30.62437979124922
$ws.Range("R11").Value = 275.619418121243
$ws.Range("S11").Value = 0.09639237773929089
$ws.Range("T11").Value = 0.09639237773929089
$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 13.68376133333333
$ws.Range("N12").Value = 41.051284
$ws.Range("O12").Value = 0.4982080363333638
$ws.Range("P12").Value = 0.4982080363333638
$ws.Range("Q12").Value = 39.65206010736756
$ws.Range("R12").Value = 356.868540966308
$ws.Range("S12").Value = 0.1248076330709106
$ws.Range("T12").Value = 0.1248076330709106
$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.6480206666666667
$ws.Range("N13").Value = 1.944062
$ws.Range("O13").Value = 0.0235935936018545
$ws.Range("P13").Value = 0.0235935936018545
$ws.Range("Q13").Value = 1.877799078743778
$ws.Range("R13").Value = 16.900191708694
$ws.Range("S13").Value = 0.00591050396287484
$ws.Range("T13").Value = 0.00591050396287484
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("M14").Value = 2.565830333333333
$ws.Range("N14").Value = 7.697490999999999
$ws.Range("O14").Value = 0.0934185609347503
$ws.Range("P14").Value = 0.0934185609347503
$ws.Range("Q14").Value = 0.7162173264788887
$ws.Range("R14").Value = 6.445955938309999
$ws.Range("S14").Value = 0.002254344138492736
$ws.Range("T14").Value = 0.002254344138492736
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("O15").Value = 0.3847798091300315
$ws.Range("P15").Value = 0.3847798091300315
$ws.Range("Q15").Value = 2.950012967665555
$ws.Range("R15").Value = 26.55011670899
$ws.Range("S15").Value = 0.009285372185603541
$ws.Range("T15").Value = 0.009285372185603543
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 13.68376133333333
$ws.Range("N16").Value = 41.051284
$ws.Range("O16").Value = 0.4982080363333638
$ws.Range("P16").Value = 0.4982080363333638
$ws.Range("Q16").Value = 3.819639526048889
$ws.Range("R16").Value = 34.37675573444
$ws.Range("S16").Value = 0.01202258261334773
$ws.Range("T16").Value = 0.01202258261334773
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6480206666666667
$ws.Range("N17").Value = 1.944062
$ws.Range("O17").Value = 0.0235935936018545
$ws.Range("P17").Value = 0.0235935936018545
$ws.Range("Q17").Value = 0.1808863288244444
$ws.Range("R17").Value = 1.62797695942
$ws.Range("S17").Value = 0.0005693523739834791
$ws.Range("T17").Value = 0.0005693523739834794

Write-Output "Applied 190 cell updates"
